$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row 1 with two new columns P and Q (14, 15), matching header style
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

# Update existing I/K/M/O columns for rows 2-25 (swap 1<->2 pattern) and add P/Q columns
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P: new column
    $ws.Cells.Item($r, 17).Value = 2  # Q: new column
}
